# Applies the "feat: add 2022-Q4 data" change:
#  1. Inserts a new worksheet "2022-Q4" (positioned right after "总计",
#     before "2022-Q3") containing the per-fund holding breakdown for
#     that quarter.
#  2. Updates the "总计" (summary) sheet: the existing quarters shift
#     down one row and a new top row with the 2022-Q4 totals is added.

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("总计")
$q3      = $wb.Worksheets.Item("2022-Q3")
$q1      = $wb.Worksheets.Item("2022-Q1")

# ---------------------------------------------------------------------
# 1. New "2022-Q4" sheet, placed before the existing "2022-Q3" sheet.
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

# Copy the header-row look & feel (bold + border + centred) from the
# "2022-Q1" sheet, which already carries the styling we need.
$q1.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats

# Copy the column-A "index" styling (bold + border + centred) too.
$q1.Range("A2:A4").Copy()
$q4.Range("A2:A11").PasteSpecial(-4122)  # xlPasteFormats

# Header row
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Make sure numeric-looking text columns stay text (so leading zeros in
# fund codes and the decimal strings survive), like the other quarter
# sheets, instead of being auto-converted to numbers.
$q4.Range("B2:B11").NumberFormat = "@"
$q4.Range("D2:G11").NumberFormat = "@"

$q4Data = @(
    @(0, "002300", "长盛医疗行业量化配置股票",       "2.79", "92.63", "6.36", "0.1774", 5),
    @(1, "000940", "富国中小盘精选混合A",             "5.50", "79.83", "2.53", "0.1392", 10),
    @(2, "014313", "鹏华创新增长一年持有期混合A",     "3.47", "60.62", "3.12", "0.1083", 8),
    @(3, "000684", "长盛养老健康产业灵活配置混合",    "1.47", "88.14", "5.83", "0.0857", 3),
    @(4, "008412", "长盛竞争优势股票A",               "0.83", "88.91", "5.64", "0.0468", 3),
    @(5, "008413", "长盛竞争优势股票C",               "0.40", "88.91", "5.64", "0.0226", 3),
    @(6, "014708", "天弘臻选健康混合A",               "0.49", "91.26", "4.61", "0.0226", 7),
    @(7, "014314", "鹏华创新增长一年持有期混合C",     "0.15", "60.62", "3.12", "0.0047", 8),
    @(8, "014709", "天弘臻选健康混合C",               "0.10", "91.26", "4.61", "0.0046", 7),
    @(9, "015690", "富国中小盘精选混合C",             "0.01", "79.83", "2.53", "0.0003", 10)
)

$row = 2
foreach ($r in $q4Data) {
    $q4.Cells.Item($row, 1).Value = $r[0]
    $q4.Cells.Item($row, 2).Value = $r[1]
    $q4.Cells.Item($row, 3).Value = $r[2]
    $q4.Cells.Item($row, 4).Value = $r[3]
    $q4.Cells.Item($row, 5).Value = $r[4]
    $q4.Cells.Item($row, 6).Value = $r[5]
    $q4.Cells.Item($row, 7).Value = $r[6]
    $q4.Cells.Item($row, 8).Value = $r[7]
    $row = $row + 1
}

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: push existing quarter rows down by one
#    and write the new 2022-Q4 totals into row 2. Values are written
#    literally (rather than copied from the cells being shifted) to
#    sidestep COM `.Value` getter quirks in this host.
# ---------------------------------------------------------------------

# Row 4 (2021-Q3) -> row 5, with a new index value of 3. Grab the
# existing index-column styling (bold + border + centred) from A4.
$summary.Range("A4").Copy()
$summary.Range("A5").PasteSpecial(-4122)  # xlPasteFormats
$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2021-Q3"
$summary.Range("C5").Value = 3
$summary.Range("D5").Value = 0.6

# Row 3 (2022-Q1) -> row 4.
$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 3
$summary.Range("D4").Value = 0.15

# Row 2 (2022-Q3) -> row 3.
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 5
$summary.Range("D3").Value = 0.28

# New 2022-Q4 totals go into row 2.
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 10
$summary.Range("D2").Value = 0.61

$summary.Activate()
